# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newer scrape counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$exhibitionUpdates = @{
    5  = 13230
    9  = 484
    11 = 994
    12 = 13790
    13 = 14416
    22 = 1096
    25 = 5469
    27 = 221
    28 = 338
    30 = 80
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row => new F value
$allTypesUpdates = @{
    5  = 13230
    10 = 484
    12 = 994
    13 = 13790
    14 = 14416
    23 = 1096
    26 = 5469
    28 = 221
    29 = 338
    31 = 80
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
